# Update "want to go" counts (column F) across sheets, reflecting the
# latest generated output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 811
$ws1.Range("F7").Value = 12066
$ws1.Range("F8").Value = 6820
$ws1.Range("F13").Value = 10
$ws1.Range("F16").Value = 103
$ws1.Range("F18").Value = 951
$ws1.Range("F23").Value = 202
$ws1.Range("F24").Value = 324
$ws1.Range("F28").Value = 55
$ws1.Range("F29").Value = 321
$ws1.Range("F30").Value = 5087
$ws1.Range("F32").Value = 1280
$ws1.Range("F33").Value = 255
$ws1.Range("F34").Value = 739
$ws1.Range("F35").Value = 1239
$ws1.Range("F36").Value = 561

# 本地生活 (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 9137
$ws3.Range("F4").Value = 1885

# 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 9137
$ws4.Range("F4").Value = 1885
$ws4.Range("F5").Value = 811
$ws4.Range("F8").Value = 12066
$ws4.Range("F9").Value = 6820
$ws4.Range("F16").Value = 10
$ws4.Range("F19").Value = 103
$ws4.Range("F21").Value = 951
$ws4.Range("F25").Value = 202
$ws4.Range("F26").Value = 324
$ws4.Range("F34").Value = 321
$ws4.Range("F36").Value = 1280
$ws4.Range("F38").Value = 255
$ws4.Range("F40").Value = 739
$ws4.Range("F41").Value = 1239
$ws4.Range("F42").Value = 561
